$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5439096666666666
$ws.Range("H2").Value = 1.631729
$ws.Range("I2").Value = 0.003493229883501837
$ws.Range("J2").Value = 0.003493229883501837
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 79.24804346473178
$ws.Range("R2").Value = 713.2323911825861
$ws.Range("S2").Value = 0.001001137939292575
$ws.Range("T2").Value = 0.001001137939292575

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5439096666666666
$ws.Range("H3").Value = 1.631729
$ws.Range("I3").Value = 0.003493229883501837
$ws.Range("J3").Value = 0.003493229883501837
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 91.81179001085911
$ws.Range("R3").Value = 826.3061100977321
$ws.Range("S3").Value = 0.001159855338196963
$ws.Range("T3").Value = 0.001159855338196963

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5439096666666666
$ws.Range("H4").Value = 1.631729
$ws.Range("I4").Value = 0.003493229883501837
$ws.Range("J4").Value = 0.003493229883501837
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 69.68905407602844
$ws.Range("R4").Value = 627.201486684256
$ws.Range("S4").Value = 0.0008803795392118843
$ws.Range("T4").Value = 0.0008803795392118843

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5439096666666666
$ws.Range("H5").Value = 1.631729
$ws.Range("I5").Value = 0.003493229883501837
$ws.Range("J5").Value = 0.003493229883501837
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 35.76808655853033
$ws.Range("R5").Value = 321.912779026773
$ws.Range("S5").Value = 0.000451857066800415
$ws.Range("T5").Value = 0.000451857066800415

$ws.Range("I6").Value = 0.00653284034046588
$ws.Range("J6").Value = 0.006532840340465881
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 148.205194766742
$ws.Range("R6").Value = 1333.846752900678
$ws.Range("S6").Value = 0.001872271374715547
$ws.Range("T6").Value = 0.001872271374715547

$ws.Range("I7").Value = 0.00653284034046588
$ws.Range("J7").Value = 0.006532840340465881
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.002169095649348447
$ws.Range("T7").Value = 0.002169095649348447

$ws.Range("I8").Value = 0.00653284034046588
$ws.Range("J8").Value = 0.006532840340465881
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 130.328515139232
$ws.Range("R8").Value = 1172.956636253088
$ws.Range("S8").Value = 0.001646435866086949
$ws.Range("T8").Value = 0.00164643586608695

$ws.Range("I9").Value = 0.00653284034046588
$ws.Range("J9").Value = 0.006532840340465881
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 66.89144618693099
$ws.Range("R9").Value = 602.023015682379
$ws.Range("S9").Value = 0.0008450374503149371
$ws.Range("T9").Value = 0.0008450374503149372

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.21601
$ws.Range("H10").Value = 0.6480299999999999
$ws.Range("I10").Value = 0.001387312330298533
$ws.Range("J10").Value = 0.001387312330298533
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 31.47281785544666
$ws.Range("R10").Value = 283.25536069902
$ws.Range("S10").Value = 0.0003975950778589872
$ws.Range("T10").Value = 0.0003975950778589872

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.21601
$ws.Range("H11").Value = 0.6480299999999999
$ws.Range("I11").Value = 0.001387312330298533
$ws.Range("J11").Value = 0.001387312330298533
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 36.46242377302666
$ws.Range("R11").Value = 328.16181395724
$ws.Range("S11").Value = 0.0004606286061054121
$ws.Range("T11").Value = 0.0004606286061054121

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.21601
$ws.Range("H12").Value = 0.6480299999999999
$ws.Range("I12").Value = 0.001387312330298533
$ws.Range("J12").Value = 0.001387312330298533
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 27.67653066954666
$ws.Range("R12").Value = 249.08877602592
$ws.Range("S12").Value = 0.0003496367060924194
$ws.Range("T12").Value = 0.0003496367060924194

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.21601
$ws.Range("H13").Value = 0.6480299999999999
$ws.Range("I13").Value = 0.001387312330298533
$ws.Range("J13").Value = 0.001387312330298533
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 14.20505067479
$ws.Range("R13").Value = 127.84545607311
$ws.Range("S13").Value = 0.0001794519402417147
$ws.Range("T13").Value = 0.0001794519402417147

$ws.Range("G14").Value = 153.9268343333333
$ws.Range("H14").Value = 461.780503
$ws.Range("I14").Value = 0.9885866174457337
$ws.Range("J14").Value = 0.9885866174457337
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 22427.2543865493
$ws.Range("R14").Value = 201845.2894789437
$ws.Range("S14").Value = 0.2833227706187171
$ws.Range("T14").Value = 0.2833227706187171

$ws.Range("G15").Value = 153.9268343333333
$ws.Range("H15").Value = 461.780503
$ws.Range("I15").Value = 0.9885866174457337
$ws.Range("J15").Value = 0.9885866174457337
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 25982.80386788793
$ws.Range("R15").Value = 233845.2348109913
$ws.Range("S15").Value = 0.3282399108429333
$ws.Range("T15").Value = 0.3282399108429333

$ws.Range("G16").Value = 153.9268343333333
$ws.Range("H16").Value = 461.780503
$ws.Range("I16").Value = 0.9885866174457337
$ws.Range("J16").Value = 0.9885866174457337
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 19722.05338314304
$ws.Range("R16").Value = 177498.4804482874
$ws.Range("S16").Value = 0.2491480548842192
$ws.Range("T16").Value = 0.2491480548842192

$ws.Range("G17").Value = 153.9268343333333
$ws.Range("H17").Value = 461.780503
$ws.Range("I17").Value = 0.9885866174457337
$ws.Range("J17").Value = 0.9885866174457337
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 10122.39471281424
$ws.Range("R17").Value = 91101.5524153282
$ws.Range("S17").Value = 0.1278758810998642
$ws.Range("T17").Value = 0.1278758810998642

